# Atualização de bases das ligas, do dia: 25-05-2024 às 15:10
#
# Swap the betting-data content (columns B and E through AB) between the
# following pairs of rows. Columns A (row index), C (Div) and D (Date)
# stay untouched in each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($sheet, $row1, $row2)

    # Swap column B (id-like numeric key)
    $rangeB1 = $sheet.Range("B$row1")
    $rangeB2 = $sheet.Range("B$row2")
    $tmpB = $rangeB1.Value()
    $rangeB1.Value() = $rangeB2.Value()
    $rangeB2.Value() = $tmpB

    # Swap columns E through AB (teams, scores, odds, etc.)
    $rangeE1 = $sheet.Range("E$row1`:AB$row1")
    $rangeE2 = $sheet.Range("E$row2`:AB$row2")
    $tmpE = $rangeE1.Value()
    $rangeE1.Value() = $rangeE2.Value()
    $rangeE2.Value() = $tmpE
}

Swap-RowData $ws 28 29
Swap-RowData $ws 101 102
Swap-RowData $ws 149 150
Swap-RowData $ws 221 222
Swap-RowData $ws 271 272
